# tambahan detail info penyelesaian proses
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet ("Sheet2") right after the existing Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Fill in the summary data
$ws2.Range("A1").Value = "Total Process"
$ws2.Range("B1").Value = 100

$ws2.Range("A2").Value = "AWT (Average Waiting Time)"
$ws2.Range("B2").Value = 786.37

$ws2.Range("A3").Value = "Total Waiting Time"
$ws2.Range("B3").Value = 78637

$ws2.Range("A4").Value = "ATAT (Average Turn Around Time)"
$ws2.Range("B4").Value = 798.95

$ws2.Range("A5").Value = "Total Turn Around Time"
$ws2.Range("B5").Value = 79895

$ws2.Range("A6").Value = "Quantum Time"
$ws2.Range("B6").Value = 12
